$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.562.45'
$ws.Range("E2").Value = '  +4.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.487.71'
$ws.Range("E3").Value = '  +2.67%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.74'
$ws.Range("E5").Value = '  +3.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.79'
$ws.Range("E6").Value = '  +4.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.484.28'
$ws.Range("E8").Value = '  +2.52%  '

$ws.Range("E9").Value = '  +7.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.31'
$ws.Range("E10").Value = '  +0.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.127'
$ws.Range("E11").Value = '  +6.49%  '

$ws.Range("E12").Value = '  +3.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.090.61'
$ws.Range("E13").Value = '  +2.84%  '

$ws.Range("E14").Value = '  -0.19%  '

$ws.Range("E15").Value = '  +4.39%  '

$ws.Range("E16").Value = '  +3.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.546.59'
$ws.Range("E17").Value = '  +4.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.495.81'
$ws.Range("E18").Value = '  +2.77%  '

$ws.Range("E19").Value = '  +3.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.09'
$ws.Range("E20").Value = '  +3.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '392.51'
$ws.Range("E21").Value = '  +4.46%  '

$ws.Range("E22").Value = '  +2.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.24'
$ws.Range("E23").Value = '  +4.21%  '

$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("E25").Value = '  +4.31%  '

$ws.Range("E26").Value = '  +5.55%  '

$ws.Range("E27").Value = '  +6.91%  '

$ws.Range("E28").Value = '  +1.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.33'
$ws.Range("E30").Value = '  +4.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.47'
$ws.Range("E31").Value = '  +5.73%  '

$ws.Range("E32").Value = '  +3.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.62'
$ws.Range("E33").Value = '  +3.32%  '

$ws.Range("E34").Value = '  +4.70%  '

$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  +9.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.89'
$ws.Range("E37").Value = '  +1.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.884'
$ws.Range("E38").Value = '  +2.93%  '

$ws.Range("E39").Value = '  +6.52%  '

$ws.Range("E40").Value = '  +6.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0744'
$ws.Range("E41").Value = '  +3.24%  '

$ws.Range("E42").Value = '  +3.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.45'
$ws.Range("E43").Value = '  +3.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.08'
$ws.Range("E44").Value = '  +4.22%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.782.57'
$ws.Range("E45").Value = '  +1.68%  '

$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.11'
$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0312'
$ws.Range("E47").Value = '  +2.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '348.18'
$ws.Range("E48").Value = '  +6.61%  '

$ws.Range("E49").Value = '  +2.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.09'
$ws.Range("E50").Value = '  +5.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.99'
$ws.Range("E51").Value = '  +14.09%  '
